$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new column before column I (9)
$ws.Columns.Item(9).EntireColumn.Insert()

# Set header text for new column I1
$ws.Range("I1").Value = "Last Version Date"
